$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = '[Noelle%Breslin%nb2565@cumc.columbia.edu%0,            Caitlin%Baptiste%NULL%1,            Cynthia%Gyamfi-Bannerman%NULL%1,            Russell%Miller%NULL%1,            Rebecca%Martinez%NULL%1,            Kyra%Bernstein%NULL%1,            Laurence%Ring%NULL%1,            Ruth%Landau%NULL%1,            Stephanie%Purisch%NULL%1,            Alexander M.%Friedman%NULL%1,            Karin%Fuchs%NULL%1,            Desmond%Sutton%NULL%1,            Maria%Andrikopoulou%NULL%1,            Devon%Rupley%NULL%1,            Jean-Ju%Sheen%NULL%1,            Janice%Aubey%NULL%1,            Noelia%Zork%NULL%1,            Leslie%Moroz%NULL%1,            Mirella%Mourad%NULL%1,            Ronald%Wapner%NULL%1,            Lynn L.%Simpson%NULL%1,            Mary E.%D’Alton%NULL%1,            Dena%Goffman%NULL%1]'
$ws.Range("E3").Value = '[Alexander%Juusela%alex.l.juusela@gmail.com%1,            Munir%Nazir%NULL%1,            Martin%Gimovsky%NULL%1]'
$ws.Range("E4").Value = '[William S.%Vintzileos%NULL%1,            Jolene%Muscat%NULL%2,            Jolene%Muscat%NULL%0,            Eva%Hoffmann%NULL%1,            Nicole S.%John%NULL%1,            Rosanne%Vertichio%NULL%1,            Anthony M.%Vintzileos%NULL%1,            Duc%Vo%NULL%2,            Duc%Vo%NULL%0]'
$ws.Range("E5").Value = '[William T.%Schnettler%william_schnettler@trihealth.com%2,            Yousef%Al Ahwel%NULL%1,            Anju%Suhag%NULL%1]'
$ws.Range("E6").Value = '[Sara N.%Iqbal%NULL%2,            Rachael%Overcash%NULL%2,            Neggin%Mokhtari%NULL%2,            Haleema%Saeed%NULL%2,            Stacey%Gold%NULL%2,            Tamika%Auguste%NULL%2,            Muhammad-Usman%Mirza%NULL%2,            Maria-Elena%Ruiz%NULL%2,            Joeffrey J.%Chahine%NULL%4,            Joeffrey J.%Chahine%NULL%0,            Masashi%Waga%NULL%4,            Masashi%Waga%NULL%0,            Glenn%Wortmann%NULL%2]'
$ws.Range("E7").Value = '[Paul C.%Browne%NULL%1,            Jennifer B.%Linfert%NULL%1,            Emilio%Perez-Jorge%NULL%1]'
$ws.Range("E8").Value = '[Jie%Yan%NULL%0,            Juanjuan%Guo%NULL%2,            Cuifang%Fan%NULL%1,            Juan%Juan%NULL%1,            Xuechen%Yu%NULL%2,            Jiafu%Li%NULL%2,            Ling%Feng%NULL%1,            Chunyan%Li%NULL%1,            Huijun%Chen%NULL%1,            Yuan%Qiao%NULL%1,            Di%Lei%NULL%1,            Chen%Wang%NULL%3,            Guoping%Xiong%NULL%1,            Fengyi%Xiao%NULL%1,            Wencong%He%NULL%1,            Qiumei%Pang%NULL%3,            Xiaoling%Hu%NULL%1,            Suqing%Wang%NULL%1,            Dunjin%Chen%NULL%1,            Yuanzhen%Zhang%NULL%2,            Liona C.%Poon%NULL%1,            Huixia%Yang%NULL%2]'
$ws.Range("E10").Value = '[Xiaotong%Wang%NULL%3,            Zhiqiang%Zhou%NULL%1,            Jianping%Zhang%NULL%1,            Fengfeng%Zhu%NULL%1,            Yongyan%Tang%NULL%1,            Xinghua%Shen%luckydrx@163.com%1]'
$ws.Range("E11").Value = '[Rong%Wen%wenxueqian_1111@163.com%1,            Yue%Sun%101799@qdfe.com%1,            Quan-Sheng%Xing%xingquansheng@qdfe.com%1]'
$ws.Range("E12").Value = '[Xinggui%Liao%NULL%1,            Huan%Yang%NULL%1,            Junfeng%Kong%NULL%2,            Junfeng%Kong%NULL%0,            Hongbing%Yang%NULL%1]'
$ws.Range("E13").Value = '[Liusheng%Hou%NULL%1,            Mingxing%Li%NULL%1,            Kai%Guo%NULL%1,            Weijia%Wang%NULL%1,            Binfei%Li%NULL%1,            Jianwei%Li%NULL%2,            Yong%Yuan%NULL%1]'
$ws.Range("E14").Value = '[Xiali%Xiong%NULL%2,            Hong%Wei%NULL%2,            Hong%Wei%NULL%0,            Zhihong%Zhang%NULL%1,            Jing%Chang%NULL%1,            Xiaopeng%Ma%NULL%1,            Xiang%Gao%NULL%1,            Qiang%Chen%qiangchen2015@sina.com%1,            Qiumei%Pang%pqm19650930@sina.com%0,            Qiumei%Pang%pqm19650930@sina.com%0]'
$ws.Range("E15").Value = '[Zhoujie%Peng%NULL%2,            Jianhui%Wang%NULL%1,            Yunbo%Mo%NULL%1,            Wei%Duan%NULL%1,            Guangjun%Xiang%NULL%1,            Ming%Yi%NULL%1,            Lei%Bao%NULL%1,            Yuan%Shi%NULL%0]'
$ws.Range("E16").Value = '[Dawei%Lu%NULL%2,            Lin%Sang%NULL%1,            Shihua%Du%NULL%1,            Tao%Li%NULL%0,            Yange%Chang%NULL%1,            Xiu‐An%Yang%yangxiuan07@mails.ucas.edu.cn%1]'
$ws.Range("E17").Value = '[Enrico M.%Ferrazzi%enrico.ferrazzi@unimi.it%0,            Luigi%Frigerio%NULL%1,            Irene%Cetin%NULL%1,            Patrizia%Vergani%NULL%1,            Arsenio%Spinillo%NULL%1,            Federico%Prefumo%NULL%1,            Edda%Pellegrini%NULL%1,            Gianluigi%Gargantini%NULL%1]'
$ws.Range("E19").Value = '[Evangelia%Vlachodimitropoulou Koumoutsea%NULL%1,            Alexandre J.%Vivanti%NULL%1,            Nadine%Shehata%NULL%1,            Alexandra%Benachi%NULL%1,            Agnes%Le Gouez%NULL%1,            Celine%Desconclois%NULL%1,            Wendy%Whittle%NULL%1,            John%Snelgrove%NULL%1,            Ann Kinga%Malinowski%NULL%1]'
$ws.Range("E20").Value = '[Lysien I.%Zambrano%NULL%1,            Itzel Carolina%Fuentes-Barahona%NULL%2,            Itzel Carolina%Fuentes-Barahona%NULL%0,            Daysi Anabell%Bejarano-Torres%NULL%1,            Carolina%Bustillo%NULL%2,            Gloria%Gonzales%NULL%1,            Gissela%Vallecillo-Chinchilla%NULL%1,            Fredal Eduardo%Sanchez-Martínez%NULL%2,            Fredal Eduardo%Sanchez-Martínez%NULL%0,            Jorge Alberto%Valle-Reconco%NULL%1,            Manuel%Sierra%NULL%2,            Manuel%Sierra%NULL%0,            D. Katterine%Bonilla-Aldana%NULL%2,            D. Katterine%Bonilla-Aldana%NULL%0,            Jaime A.%Cardona-Ospina%NULL%2,            Jaime A.%Cardona-Ospina%NULL%0,            Alfonso J.%Rodríguez-Morales%NULL%1]'
$ws.Range("E21").Value = '[Dong Hwan%Lee%NULL%4,            Jihyang%Lee%NULL%6,            Jihyang%Lee%NULL%0,            Eunju%Kim%NULL%6,            Eunju%Kim%NULL%0,            Kyeongyoon%Woo%NULL%6,            Kyeongyoon%Woo%NULL%0,            Hak Youle%Park%NULL%6,            Hak Youle%Park%NULL%0,            Jihyun%An%NULL%6,            Jihyun%An%NULL%0]'
$ws.Range("E23").Value = '[AZAP%ALPAY%coreGivesNoEmail%2,           Kalafat%Erkan%coreGivesNoEmail%1,           KO\u00c7%ASLI%coreGivesNoEmail%1,           Ozisik%S.%coreGivesNoEmail%1,           UZUN%\u00c7A\u011eLAR%coreGivesNoEmail%1,           VARLI%BULUT%coreGivesNoEmail%1,           Yaprak%E.%coreGivesNoEmail%1,           \u00c7INAR%G\u00dcLE%coreGivesNoEmail%1]'
$ws.Range("E24").Value = '[Parisa%Karami%NULL%1,            Maliheh%Naghavi%NULL%1,            Abdolamir%Feyzi%NULL%1,            Mehdi%Aghamohammadi%NULL%1,            Mohammad Sadegh%Novin%NULL%1,            Ahmadreza%Mobaien%NULL%1,            Mohamad%Qorbanisani%NULL%1,            Aida%Karami%NULL%1,            Amir Hossein%Norooznezhad%norooznezhad@gmail.com%1]'
$ws.Range("E25").Value = '[Belinda%Lowe%Belinda.lowe2@health.qld.gov.au%1,            Benjamin%Bopp%NULL%2,            Benjamin%Bopp%NULL%0]'
$ws.Range("E26").Value = '[D.%González Romero%NULL%1,            J.%Ocampo Pérez%NULL%1,            L.%González Bautista%NULL%1,            L.%Santana-Cabrera%NULL%1]'
$ws.Range("E27").Value = '[Maria Claudia%Alzamora%NULL%3,            Tania%Paredes%NULL%4,            Tania%Paredes%NULL%0,            David%Caceres%NULL%2,            Camille M.%Webb%NULL%2,            Luis M.%Valdez%NULL%2,            Mauricio%La Rosa%NULL%2]'
$ws.Range("E28").Value = '[K. Aparna%Sharma%NULL%1,            Rajesh%Kumari%NULL%1,            Garima%Kachhawa%NULL%1,            Anjolie%Chhabra%NULL%1,            Ramesh%Agarwal%NULL%1,            Akash%Sharma%NULL%1,            Sunesh%Kumar%NULL%1,            Neerja%Bhatla%neerja.bhatla07@gmail.com%1]'
